$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.765799999999992
$ws.Range("B4").Value = 4.666300000000003
$ws.Range("E6").Value = 12.3078
$ws.Range("B7").Value = 5.709700000000002
$ws.Range("E7").Value = 12.32499999999999
$ws.Range("B8").Value = 5.638299999999991
$ws.Range("E8").Value = 13.3638
$ws.Range("A11").Value = -21.81040000000002
$ws.Range("A12").Value = -22.7628
$ws.Range("B12").Value = 6.219599999999998
$ws.Range("B14").Value = 8.652400000000005
$ws.Range("A15").Value = -21.44880000000002
$ws.Range("E19").Value = 12.8739
$ws.Range("E21").Value = 12.6899
$ws.Range("B22").Value = 5.194100000000002
$ws.Range("E24").Value = 12.76329999999999
$ws.Range("E25").Value = 12.995
